$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("table", "zhuōzi", "桌子", 4),
    @("these", "zhè xiē", "这些", 9),
    @("those", "nà xiē", "那些", 9),
    @("which ones", "nǎ xiē", "哪些", 9),
    @("ticket", "piào", "票", 9),
    @("newspaper", "bàozhǐ", "报纸", 18),
    @("movie", "diànyǐng", "电影", 18)
)

$startRow = 204
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$ws.Range("G9").Select() | Out-Null
